$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add I1 = "I0" and J1 = "IF" ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (style) of the existing H1 header cell onto the two
# new header cells so they share the same cell style (border/bold/alignment).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# --- Data rows (2-11): add values for the new I and J columns ---
$values = @(
    @(1, 4),
    @(1, 4),
    @(1, 5),
    @(1, 7),
    @(1, 5),
    @(1, 5),
    @(1, 4),
    @(6, 8),
    @(6, 7),
    @(3, 4)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
